# School name is part of default values in Connectors.xlsx
#
# "modules__school__schoolName" is already captured as a per-connector
# column on the "Connectors" sheet, so it does not belong among the
# global/default values on the "Default Values" sheet as well. Drop that
# (now redundant) column there, and tidy up the no-longer-needed explicit
# header styling on the "Connectors" sheet.

$wb = $excel.ActiveWorkbook

$wsConnectors = $wb.Worksheets.Item("Connectors")
$wsDefaults   = $wb.Worksheets.Item("Default Values")

# --- "Default Values" sheet: drop the schoolName column (column G). ---
# This shifts the "autoMailAfterOnboarding" / "autoMailBeforeOffboarding"
# columns left by one, and re-bolds them with the sheet's existing header
# style (reusing the same style as the other header cells).
$wsDefaults.Columns.Item(7).Delete()
$wsDefaults.Range("F1:H1").Font.Bold = $true

# Move the active selection.
$wsDefaults.Activate()
$wsDefaults.Range("F6").Select()

# --- "Connectors" sheet: clear the (redundant) explicit style on the
# header row cells A1:H1 so they fall back to the default style, matching
# the unstyled header cells I1:L1 already on that row. ---
$wsConnectors.Range("A1:H1").Style = "Normal"

$wb.Save()
